$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# "IMAGE: add BMP support"
#
# The document lists one "+++IMAGE (...)+++" placeholder paragraph per
# supported image extension. This change:
#   1) simplifies the (oddly split, bookmark-interrupted) ".jpeg" run
#      into a single clean run,
#   2) keeps the ".gif" paragraph as-is, and
#   3) adds a brand new ".bmp" paragraph after it - reusing the same
#      "typed in multiple passes" run layout (and the _GoBack bookmark)
#      that the ".jpeg" paragraph used to have.
# ----------------------------------------------------------------------

# Locate the ".jpeg" and ".gif" placeholder paragraphs by content instead
# of hard-coded indices, so the script keeps working even if unrelated
# paragraphs get added/removed earlier in the document.
$jpegParaIndex = -1
$gifParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*sample.jpeg*") {
        $jpegParaIndex = $i
    } elseif ($t -like "*sample.gif*") {
        $gifParaIndex = $i
    }
}

# --- 1) Simplify the ".jpeg" paragraph down to a single run --------------
# This also removes the _GoBack bookmark that used to sit inside it - it
# gets recreated below, inside the new ".bmp" paragraph.
$pJpeg = $d.Paragraphs.Item($jpegParaIndex)
$jpegRange = $d.Range($pJpeg.Range.Start, $pJpeg.Range.End)
$jpegRange.Text = "+++IMAGE ({ width: 3, height: 3, path: './sample.jpeg' })+++"

# --- 2) Build the new ".bmp" paragraph right after ".gif" ------------------
$pGif = $d.Paragraphs.Item($gifParaIndex)

# New paragraph #1: holds "+++IMAGE ({ width: 3, height: 3, path: "
$pGif.Range.InsertParagraphAfter()
$pPart1 = $d.Paragraphs.Item($gifParaIndex + 1)
$part1Start = $pPart1.Range.Start
$part1Text = "+++IMAGE ({ width: 3, height: 3, path: "
$d.Range($part1Start, $part1Start).InsertAfter($part1Text)

# New paragraph #2 (temporary): holds "'./sample.bmp" then "' })+++",
# split into two runs via a bookmark (mirrors the original ".jpeg" layout)
$pPart1 = $d.Paragraphs.Item($gifParaIndex + 1)
$pPart1.Range.InsertParagraphAfter()
$pPart2 = $d.Paragraphs.Item($gifParaIndex + 2)
$part2Start = $pPart2.Range.Start
$part2Text = "'./sample.bmp"
$d.Range($part2Start, $part2Start).InsertAfter($part2Text)

$afterPart2 = $part2Start + $part2Text.Length
$part3Text = "' })+++"
$d.Range($afterPart2, $afterPart2).InsertAfter($part3Text)

# Recreate the _GoBack bookmark between "'./sample.bmp" and "' })+++"
$d.Bookmarks.Add("_GoBack", $d.Range($afterPart2, $afterPart2)) | Out-Null

# Merge the two helper paragraphs into one (deletes the paragraph mark
# between them) - this keeps "...path: " and "'./sample.bmp" as two
# distinct runs instead of letting them collapse into one.
$pPart1 = $d.Paragraphs.Item($gifParaIndex + 1)
$markPos = $pPart1.Range.End - 1
$d.Range($markPos, $pPart1.Range.End).Delete()
